$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.294.40'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '2.270.76'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'113.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.65%  '
$ws.Range("D6").Value = "'264.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.28%  '
$ws.Range("E7").Value = '  -0.94%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = "'0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.26%  '
$ws.Range("D10").Value = "'48.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").Value = "'0.0927"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").Value = "'8.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.31%  '
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("D15").Value = '2.609.51'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '2.271.21'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").Value = '43.158.23'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("D20").Value = "'6.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.96%  '
$ws.Range("D21").Value = "'71.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").Value = "'2.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("D23").Value = "'9.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.89%  '
$ws.Range("D24").Value = "'230.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("D25").Value = "'2.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = "'11.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").Value = "'40.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.79%  '
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").Value = "'171.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.73%  '
$ws.Range("E33").Value = '  -3.24%  '
$ws.Range("D34").Value = "'0.0907"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("D35").Value = "'5.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = "'4.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").Value = "'0.0351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.37%  '
$ws.Range("D39").Value = "'3.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").Value = '  -7.67%  '
$ws.Range("D41").Value = "'14.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.48%  '
$ws.Range("D42").Value = "'75.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.36%  '
$ws.Range("D43").Value = "'2.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.67%  '
$ws.Range("E44").Value = '  -1.41%  '
$ws.Range("D45").Value = "'6.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.05%  '
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("D48").Value = "'8.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("D49").Value = "'0.0985"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.70%  '
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("E51").Value = '  +0.52%  '
